$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bollinger Bands")

# Rename the header of the "testdata" table's first column from "index" to "i".
# Setting the header cell's value updates the table's column name as well.
$table = $ws.ListObjects.Item("testdata")
$ws.Range("A1").Value2 = "i"

# Decrement the data values in column A (the "i"/"index" column) by 1,
# turning the 1-based index (1..502) into a 0-based index (0..501).
$lastRow = $table.ListRows.Count + 1
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Narrow column A now that values are shorter (max 3 digits instead of "index").
# (ColumnWidth uses character units; stored width = ColumnWidth + 5/6, so subtract
# 5/6 to land exactly on a stored/raw width of 4.)
$ws.Columns.Item(1).ColumnWidth = 4 - 5/6
